# Workbook: flotta (fleet) tracker.
# Vehicle GL350TJ's operator "VALENTINA.DE.GREGORIIS" is being reassigned
# ("DA ASSEGNARE (MATERNITA')" = "to be assigned (maternity leave)"),
# effective 2025-12-18. Update the current-state sheet and the matching
# entry on the history sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Stato Attuale")
$ws2 = $wb.Worksheets.Item("Storico Passaggi")

# --- Sheet1 "Stato Attuale": row 42 is vehicle GL350TJ ---
# Update current operator and stamp the assignment date in column C
# (leading apostrophe keeps the ISO date as literal text instead of
# letting it be auto-converted to a date serial number).
$ws1.Cells.Item(42, 2).Value = "DA ASSEGNARE (MATERNITA')"
$ws1.Cells.Item(42, 3).Value = "'2025-12-18"
$ws1.Cells.Item(42, 3).Style = "Normal"

# --- Sheet2 "Storico Passaggi": log this change on row 2 ---
$ws2.Cells.Item(2, 1).Value = "GL350TJ"
$ws2.Cells.Item(2, 2).Value = "VALENTINA.DE.GREGORIIS"
$ws2.Cells.Item(2, 3).Value = "DA ASSEGNARE (MATERNITA')"
# Data_Cambio (D2) already reads 2025-12-18 -- leave it untouched.
